# Add an "InvalidLogin" sheet (mirroring "ValidLogin") with a bad
# username/password pair, make it the active sheet/tab, and drop the
# tab-selected state + stale cell selection from "ValidLogin".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ValidLogin keeps zoom etc. but its selection becomes the whole used range
# (A1:B2) and it is no longer the selected/active tab.
$ws1.Range("A1:B2").Select() | Out-Null

# Insert the new sheet immediately after "ValidLogin".
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "InvalidLogin"

# Same shape as ValidLogin: UserName/Password headers + one data row,
# this time with an invalid credential pair.
$ws2.Range("A1").Value = "UserName"
$ws2.Range("B1").Value = "Password"
$ws2.Range("A2").Value = "Bhanu"
$ws2.Range("B2").Value = "Damager"

# Match the column widths used on the original sheet (best-fit to content).
$ws2.Columns.Item(1).ColumnWidth = 10.28515625
$ws2.Columns.Item(2).ColumnWidth = 9.42578125

# InvalidLogin is left as the active tab/sheet, zoomed in further, with
# the cursor sitting just below the data (B3).
$excel.ActiveWindow.Zoom = 250
$ws2.Range("B3").Select() | Out-Null
